$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new student row (row 3)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Igor Kim"
$ws.Range("C3").Value = "i@yandex.ru"
$ws.Range("D3").Value = "23CST4"

# Remove the frozen header pane
$ws.Application.ActiveWindow.FreezePanes = $false

$ws.Range("A1").Select()
